$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values for the rows that were repulled/recalculated.
$ws.Range("F3").Value = -3
$ws.Range("F9").Value = -5
$ws.Range("F13").Value = -1
$ws.Range("F15").Value = -4
$ws.Range("F16").Value = -1
$ws.Range("F18").Value = -2
